# Apply the two changes described by the diff:
#   1. Slide 6's table switches to a different built-in table style.
#   2. The presentation's theme ("Integral" green palette) is swapped out
#      for the stock Office default ("Office Theme") palette.

$p = $ppt.ActivePresentation

# --- 1. Change the table's style GUID on slide 6 ---------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{271BD433-C645-470B-98F3-C17CE0C7C46D}")

# --- 2. Swap the theme color palette from "Integral" to "Office Theme" -----
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# VBA-style RGB() packs as 0x00BBGGRR (Blue<<16 | Green<<8 | Red)
function RGBVal($r, $g, $b) { return ($b * 65536) + ($g * 256) + $r }

$colorScheme.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      -> 000000
$colorScheme.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$colorScheme.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      -> 44546A
$colorScheme.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$colorScheme.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$colorScheme.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  -> ED7D31
$colorScheme.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$colorScheme.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  -> FFC000
$colorScheme.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  -> 4472C4
$colorScheme.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  -> 70AD47
$colorScheme.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    -> 0563C1
$colorScheme.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink -> 954F72
